$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO")
[void]$ws.Activate()

# --- Row 14 / B14: status changes from "offen" to a new "verschoben" status ---
# Start from the same base formatting family as the other status pills
# (white font on a solid theme fill) by copying an existing "done" cell's
# format, then re-point the fill to a new (7th) theme color.
$ws.Range("B9").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Interior.ThemeColor = 6
$ws.Range("B14").Interior.TintAndShade = -0.249977111117893
$ws.Range("B14").Value = "verschoben"

# --- New row 23: note about price now being included in confirmation emails ---
$ws.Range("A23").Value = "Preis in Email: Preis an sendMail in mailController übergeben, in mailController eigene price-Hashlist anlegen, die nach userID sortiert ist. Beim Senden der Mail den Preis rausholen und Eintrag löschen. Preis wir din requestController berechnet"

$ws.Range("B9").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "done"

$excel.CutCopyMode = 0

[void]$ws.Range("B23").Select()
